# Applies the authored changes to template_opp_analysis.pptx:
#   - Slide 1, "Rectangle 53" (placeholder {Corners_right_positions_vs}):
#       nudge vertical position, fix placeholder typo -> {Corners_right_positions_vIs}
#   - Slide 1, "Rectangle 42" (placeholder {Corners_right_shots_vs}):
#       nudge position, fix placeholder typo -> {Corners_right_shots_vIs}
#   - Slide 1, "Rectangle 46" (placeholder {LOGO):
#       nudge vertical position, close the unterminated placeholder -> {LOGO}
#   - Slide 2, "Rectangle 40" (placeholder Def_corners_headers}):
#       nudge vertical position, fix the placeholder braces -> {def_corners_headers}
#
# NOTE: the shapes use <a:spAutoFit/>, so re-setting the run text makes
# PowerPoint recompute the text-box height; each shape's Height is
# re-asserted right after the text edit so only the intended <a:off>/<a:t>
# values move, matching the authored OOXML diff exactly.

$p = $ppt.ActivePresentation

# ---- Slide 1 ----
$s1 = $p.Slides.Item(1)

# Rectangle 53 (id 54): {Corners_right_positions_vs} -> {Corners_right_positions_vIs}
$sh = $s1.Shapes.Item("Rectangle 53")
$sh.TextFrame.TextRange.Runs(1).Text = "{Corners_right_positions_vIs}"
$sh.Height = 126.0
$sh.Top = 135.9792938232422

# Rectangle 42 (id 43): {Corners_right_shots_vs} -> {Corners_right_shots_vIs}
$sh = $s1.Shapes.Item("Rectangle 42")
$sh.TextFrame.TextRange.Runs(1).Text = "{Corners_right_shots_vIs}"
$sh.Height = 126.0
$sh.Left = 304.1452941894531
$sh.Top = 339.1224670410156

# Rectangle 46 (id 47): {LOGO -> {LOGO}
$sh = $s1.Shapes.Item("Rectangle 46")
$sh.TextFrame.TextRange.Runs(1).Text = "{LOGO}"
$sh.Height = 63.36000061035156
$sh.Top = 5.620630264282227

# ---- Slide 2 ----
$s2 = $p.Slides.Item(2)

# Rectangle 40 (id 41): Def_corners_headers} -> {def_corners_headers}
$sh = $s2.Shapes.Item("Rectangle 40")
$sh.TextFrame.TextRange.Runs(1).Text = "{def_corners_headers}"
$sh.Height = 200.8800048828125
$sh.Top = 596.3394165039062
